$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'Datos actualizados a 14 de Mayo de 2020 a las 23:35'
$ws.Range("B4").Value = 1451544
$ws.Range("C4").Value = 21196
$ws.Range("E4").Value = 1048764
$ws.Range("B11").Value = 174950
$ws.Range("C11").Value = 852
$ws.Range("E11").Value = 16724
$ws.Range("G11").Value = 65
$ws.Range("H11").Value = 7926
$ws.Range("B82").Value = 1971
$ws.Range("C82").Value = 59
$ws.Range("D82").Value = 930
$ws.Range("E82").Value = 1017
$ws.Range("A93").Value = 'Somalia'
$ws.Range("C93").Value = 65
$ws.Range("D93").Value = 135
$ws.Range("E93").Value = 1096
$ws.Range("F93").Value = 2
$ws.Range("G93").Value = 1
$ws.Range("H93").Value = 53
$ws.Range("A94").Value = 'Republica de Yibuti'
$ws.Range("B94").Value = 1284
$ws.Range("C94").Value = 16
$ws.Range("D94").Value = 905
$ws.Range("E94").Value = 376
$ws.Range("H94").Value = 3
$ws.Range("A95").Value = 'Consejo Danes para los Refugiados'
$ws.Range("B95").Value = 1242
$ws.Range("C95").Value = 73
$ws.Range("D95").Value = 157
$ws.Range("E95").Value = 1035
$ws.Range("H95").Value = 50
$ws.Range("B96").Value = 1210
$ws.Range("C96").Value = 67
$ws.Range("E96").Value = 567
$ws.Range("F96").Value = 10
$ws.Range("G96").Value = 2
$ws.Range("H96").Value = 16
$ws.Range("A159").Value = 'Bermudas'
$ws.Range("C159").Value = 1
$ws.Range("D159").Value = 66
$ws.Range("E159").Value = 47
$ws.Range("F159").Value = 2
$ws.Range("G159").Value = 1
$ws.Range("H159").Value = 9
$ws.Range("A160").Value = 'Camboya'
$ws.Range("B160").Value = 122
$ws.Range("D160").Value = 121
$ws.Range("E160").Value = 1
$ws.Range("F160").Value = 1
$ws.Range("H160").Value = 0
$ws.Range("D164").Value = 93
$ws.Range("E164").Value = 5
